$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1031194.75
$ws.Cells.Item(17, 9).Value = 2457.3333
$ws.Cells.Item(17, 10).Value = 3088669.8
$ws.Cells.Item(17, 11).Value = 7371.999899999999
$ws.Cells.Item(17, 12).Value = 9266009.399999999
$ws.Cells.Item(17, 13).Value = -7203.999899999999
$ws.Cells.Item(17, 14).Value = -9266345.399999999

$ws.Cells.Item(96, 8).Value = 2596.875
$ws.Cells.Item(96, 9).Value = 2689.2
$ws.Cells.Item(96, 11).Value = 8067.599999999999
$ws.Cells.Item(96, 13).Value = -6694.599999999999

$ws.Cells.Item(116, 8).Value = 19535.5
$ws.Cells.Item(116, 9).Value = 18992.6
$ws.Cells.Item(116, 11).Value = 18992.6
$ws.Cells.Item(116, 13).Value = -15550.6

$ws.Cells.Item(132, 8).Value = 27936.63
$ws.Cells.Item(132, 9).Value = 1510.9667
$ws.Cells.Item(132, 11).Value = 4532.9001
$ws.Cells.Item(132, 13).Value = -2002.9001

$ws.Cells.Item(133, 8).Value = 54490
$ws.Cells.Item(133, 10).Value = 54490
$ws.Cells.Item(133, 12).Value = 54490
$ws.Cells.Item(133, 14).Value = -64610

$ws.Cells.Item(135, 8).Value = 1190.7391
$ws.Cells.Item(135, 9).Value = 1077.1765
$ws.Cells.Item(135, 10).Value = 1512.5
$ws.Cells.Item(135, 11).Value = 9694.5885
$ws.Cells.Item(135, 12).Value = 13612.5
$ws.Cells.Item(135, 13).Value = -7159.5885
$ws.Cells.Item(135, 14).Value = -18682.5

$ws.Cells.Item(137, 8).Value = 1606.4706
$ws.Cells.Item(137, 9).Value = 1433.25
$ws.Cells.Item(137, 11).Value = 4299.75
$ws.Cells.Item(137, 13).Value = -1749.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2453.7856
$ws.Cells.Item(122, 9).Value = 2362.75
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 7088.25
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -4638.25
$ws.Cells.Item(122, 14).Value = -13900

$ws.Cells.Item(135, 8).Value = 112266.336
$ws.Cells.Item(135, 10).Value = 112266.336
$ws.Cells.Item(135, 12).Value = 112266.336
$ws.Cells.Item(135, 14).Value = -122406.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 401526.16
$ws.Cells.Item(22, 9).Value = 769.4375
$ws.Cells.Item(22, 10).Value = 532385.5
$ws.Cells.Item(22, 11).Value = 769.4375
$ws.Cells.Item(22, 12).Value = 532385.5
$ws.Cells.Item(22, 13).Value = -596.4375
$ws.Cells.Item(22, 14).Value = -532731.5

$ws.Cells.Item(81, 8).Value = 22481.555
$ws.Cells.Item(81, 10).Value = 22481.555
$ws.Cells.Item(81, 12).Value = 22481.555
$ws.Cells.Item(81, 14).Value = -24603.555

$ws.Cells.Item(84, 8).Value = 22481.555
$ws.Cells.Item(84, 10).Value = 22481.555
$ws.Cells.Item(84, 12).Value = 67444.66500000001
$ws.Cells.Item(84, 14).Value = -78052.66500000001

$ws.Cells.Item(138, 8).Value = 88950
$ws.Cells.Item(138, 10).Value = 88950
$ws.Cells.Item(138, 12).Value = 88950
$ws.Cells.Item(138, 14).Value = -99230

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2182.0833
$ws.Cells.Item(58, 9).Value = 2107.9092
$ws.Cells.Item(58, 11).Value = 2107.9092
$ws.Cells.Item(58, 13).Value = -1904.9092

$ws.Cells.Item(62, 8).Value = 7027.7334
$ws.Cells.Item(62, 9).Value = 6454.8887
$ws.Cells.Item(62, 11).Value = 6454.8887
$ws.Cells.Item(62, 13).Value = -5830.8887

$ws.Cells.Item(65, 8).Value = 7027.7334
$ws.Cells.Item(65, 9).Value = 6454.8887
$ws.Cells.Item(65, 11).Value = 32274.4435
$ws.Cells.Item(65, 13).Value = -29154.4435

$ws.Cells.Item(107, 8).Value = 1267
$ws.Cells.Item(107, 9).Value = 1312.875
$ws.Cells.Item(107, 10).Value = 900
$ws.Cells.Item(107, 11).Value = 1312.875
$ws.Cells.Item(107, 12).Value = 900
$ws.Cells.Item(107, 13).Value = 607.125
$ws.Cells.Item(107, 14).Value = -4740

$ws.Cells.Item(122, 8).Value = 1404.1428
$ws.Cells.Item(122, 9).Value = 1101.6
$ws.Cells.Item(122, 10).Value = 2160.5
$ws.Cells.Item(122, 11).Value = 3304.8
$ws.Cells.Item(122, 12).Value = 6481.5
$ws.Cells.Item(122, 13).Value = -854.7999999999997
$ws.Cells.Item(122, 14).Value = -11381.5

$ws.Cells.Item(132, 8).Value = 2139.2778
$ws.Cells.Item(132, 9).Value = 2111.8484
$ws.Cells.Item(132, 10).Value = 2441
$ws.Cells.Item(132, 11).Value = 6335.5452
$ws.Cells.Item(132, 12).Value = 7323
$ws.Cells.Item(132, 13).Value = -3805.5452
$ws.Cells.Item(132, 14).Value = -12383

$ws.Cells.Item(134, 8).Value = 25823.926
$ws.Cells.Item(134, 9).Value = 7557.9
$ws.Cells.Item(134, 11).Value = 22673.7
$ws.Cells.Item(134, 13).Value = -20138.7

$ws.Cells.Item(136, 8).Value = 2182.0833
$ws.Cells.Item(136, 9).Value = 2107.9092
$ws.Cells.Item(136, 11).Value = 6323.7276
$ws.Cells.Item(136, 13).Value = -3773.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 8396
$ws.Cells.Item(3, 9).Value = 3993.3333
$ws.Cells.Item(3, 11).Value = 11979.9999
$ws.Cells.Item(3, 13).Value = -11867.9999

$ws.Cells.Item(34, 8).Value = 3079.4
$ws.Cells.Item(34, 10).Value = 4356.5713
$ws.Cells.Item(34, 12).Value = 13069.7139
$ws.Cells.Item(34, 14).Value = -13237.7139

$ws.Cells.Item(39, 8).Value = 6563.2812
$ws.Cells.Item(39, 9).Value = 1952.7
$ws.Cells.Item(39, 10).Value = 8659
$ws.Cells.Item(39, 11).Value = 5858.1
$ws.Cells.Item(39, 12).Value = 25977
$ws.Cells.Item(39, 13).Value = -5564.1
$ws.Cells.Item(39, 14).Value = -26565

$ws.Cells.Item(55, 8).Value = 2695.3
$ws.Cells.Item(55, 10).Value = 2950.3333
$ws.Cells.Item(55, 12).Value = 8850.999899999999
$ws.Cells.Item(55, 14).Value = -9204.999899999999

$ws.Cells.Item(107, 8).Value = 858.4
$ws.Cells.Item(107, 9).Value = 673
$ws.Cells.Item(107, 11).Value = 2019
$ws.Cells.Item(107, 13).Value = -99

$ws.Cells.Item(132, 8).Value = 1401.1724
$ws.Cells.Item(132, 9).Value = 1329.68
$ws.Cells.Item(132, 11).Value = 11967.12
$ws.Cells.Item(132, 13).Value = -9437.120000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 8500
$ws.Cells.Item(20, 9).Value = 4750
$ws.Cells.Item(20, 11).Value = 4750
$ws.Cells.Item(20, 13).Value = -4505

$ws.Cells.Item(24, 8).Value = 19142.143
$ws.Cells.Item(24, 9).Value = 4500
$ws.Cells.Item(24, 11).Value = 4500
$ws.Cells.Item(24, 13).Value = -4327

$ws.Cells.Item(122, 8).Value = 2966.2942
$ws.Cells.Item(122, 9).Value = 3027.25
$ws.Cells.Item(122, 10).Value = 2820
$ws.Cells.Item(122, 11).Value = 9081.75
$ws.Cells.Item(122, 12).Value = 8460
$ws.Cells.Item(122, 13).Value = -6631.75
$ws.Cells.Item(122, 14).Value = -13360

$ws.Cells.Item(134, 8).Value = 44367.89
$ws.Cells.Item(134, 10).Value = 44367.89
$ws.Cells.Item(134, 12).Value = 133103.67
$ws.Cells.Item(134, 14).Value = -138173.67

$ws.Cells.Item(136, 8).Value = 32911.09
$ws.Cells.Item(136, 10).Value = 32911.09
$ws.Cells.Item(136, 12).Value = 98733.26999999999
$ws.Cells.Item(136, 14).Value = -103833.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 6999.6665
$ws.Cells.Item(42, 10).Value = 7000
$ws.Cells.Item(42, 12).Value = 7000
$ws.Cells.Item(42, 14).Value = -8126

$ws.Cells.Item(46, 8).Value = 400
$ws.Cells.Item(46, 10).Value = 400
$ws.Cells.Item(46, 12).Value = 400
$ws.Cells.Item(46, 14).Value = -776

$ws.Cells.Item(49, 8).Value = 6999.6665
$ws.Cells.Item(49, 10).Value = 7000
$ws.Cells.Item(49, 12).Value = 7000
$ws.Cells.Item(49, 14).Value = -7294

$ws.Cells.Item(55, 8).Value = 193.47058
$ws.Cells.Item(55, 9).Value = 170.41667
$ws.Cells.Item(55, 11).Value = 170.41667
$ws.Cells.Item(55, 13).Value = 2.583329999999989

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5449.5835
$ws.Cells.Item(62, 9).Value = 4123.8335
$ws.Cells.Item(62, 10).Value = 6775.3335
$ws.Cells.Item(62, 11).Value = 4123.8335
$ws.Cells.Item(62, 12).Value = 6775.3335
$ws.Cells.Item(62, 13).Value = -3499.8335
$ws.Cells.Item(62, 14).Value = -8023.3335

$ws.Cells.Item(65, 8).Value = 5449.5835
$ws.Cells.Item(65, 9).Value = 4123.8335
$ws.Cells.Item(65, 10).Value = 6775.3335
$ws.Cells.Item(65, 11).Value = 20619.1675
$ws.Cells.Item(65, 12).Value = 33876.6675
$ws.Cells.Item(65, 13).Value = -17499.1675
$ws.Cells.Item(65, 14).Value = -40116.6675

$ws.Cells.Item(122, 8).Value = 2337.652
$ws.Cells.Item(122, 9).Value = 2085.5
$ws.Cells.Item(122, 10).Value = 3245.4
$ws.Cells.Item(122, 11).Value = 6256.5
$ws.Cells.Item(122, 12).Value = 9736.200000000001
$ws.Cells.Item(122, 13).Value = -3806.5
$ws.Cells.Item(122, 14).Value = -14636.2

$ws.Cells.Item(135, 8).Value = 81680
$ws.Cells.Item(135, 10).Value = 81680
$ws.Cells.Item(135, 12).Value = 81680
$ws.Cells.Item(135, 14).Value = -91820
